$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "self-pair" rows (Sending cluster == Target cluster): original rows 2 (ECs/ECs), 6 (FAPs/FAPs), 10 (MuSCs/MuSCs)
# Delete bottom-to-top so earlier row indices stay valid as later rows shift up.
$ws.Rows(10).Delete()
$ws.Rows(6).Delete()
$ws.Rows(2).Delete()

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntn1"
$ws.Range("C2").Value = "Adora2b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.644726333333334
$ws.Range("H2").Value = 4.934179
$ws.Range("I2").Value = 0.03084360558270512
$ws.Range("J2").Value = 0.03084360558270512
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.467117666666667
$ws.Range("N2").Value = 4.401353
$ws.Range("O2").Value = 0.2087313363221112
$ws.Range("P2").Value = 0.2087313363221112
$ws.Range("Q2").Value = 2.413007060465223
$ws.Range("R2").Value = 21.717063544187
$ws.Range("S2").Value = 0.006438027010270168
$ws.Range("T2").Value = 0.006438027010270168

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntn1"
$ws.Range("C3").Value = "Adora2b"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.644726333333334
$ws.Range("H3").Value = 4.934179
$ws.Range("I3").Value = 0.03084360558270512
$ws.Range("J3").Value = 0.03084360558270512
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.561619333333333
$ws.Range("N3").Value = 16.684858
$ws.Range("O3").Value = 0.7912686636778888
$ws.Range("P3").Value = 0.7912686636778888
$ws.Range("Q3").Value = 9.147341773509112
$ws.Range("R3").Value = 82.32607596158199
$ws.Range("S3").Value = 0.02440557857243495
$ws.Range("T3").Value = 0.02440557857243495

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ntn1"
$ws.Range("C4").Value = "Adora2b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 44.154177
$ws.Range("H4").Value = 132.462531
$ws.Range("I4").Value = 0.828024694817689
$ws.Range("J4").Value = 0.828024694817689
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.467117666666667
$ws.Range("N4").Value = 4.401353
$ws.Range("O4").Value = 0.2087313363221112
$ws.Range("P4").Value = 0.2087313363221112
$ws.Range("Q4").Value = 64.77937313382701
$ws.Range("R4").Value = 583.0143582044431
$ws.Range("S4").Value = 0.1728347010570045
$ws.Range("T4").Value = 0.1728347010570045

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntn1"
$ws.Range("C5").Value = "Adora2b"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 44.154177
$ws.Range("H5").Value = 132.462531
$ws.Range("I5").Value = 0.828024694817689
$ws.Range("J5").Value = 0.828024694817689
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.561619333333333
$ws.Range("N5").Value = 16.684858
$ws.Range("O5").Value = 0.7912686636778888
$ws.Range("P5").Value = 0.7912686636778888
$ws.Range("Q5").Value = 245.568724450622
$ws.Range("R5").Value = 2210.118520055598
$ws.Range("S5").Value = 0.6551899937606845
$ws.Range("T5").Value = 0.6551899937606845

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ntn1"
$ws.Range("C6").Value = "Adora2b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7.525807
$ws.Range("H6").Value = 22.577421
$ws.Range("I6").Value = 0.1411316995996059
$ws.Range("J6").Value = 0.1411316995996059
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.467117666666667
$ws.Range("N6").Value = 4.401353
$ws.Range("O6").Value = 0.2087313363221112
$ws.Range("P6").Value = 0.2087313363221112
$ws.Range("Q6").Value = 11.04124440562367
$ws.Range("R6").Value = 99.37119965061301
$ws.Range("S6").Value = 0.0294586082548365
$ws.Range("T6").Value = 0.0294586082548365

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ntn1"
$ws.Range("C7").Value = "Adora2b"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7.525807
$ws.Range("H7").Value = 22.577421
$ws.Range("I7").Value = 0.1411316995996059
$ws.Range("J7").Value = 0.1411316995996059
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.561619333333333
$ws.Range("N7").Value = 16.684858
$ws.Range("O7").Value = 0.7912686636778888
$ws.Range("P7").Value = 0.7912686636778888
$ws.Range("Q7").Value = 41.85567371013533
$ws.Range("R7").Value = 376.701063391218
$ws.Range("S7").Value = 0.1116730913447694
$ws.Range("T7").Value = 0.1116730913447694

